# Upd: Add paras allowing bots use Semi/Burst weapon correctly Updated Docs
#
# Adds a new "BotReload" row (row 13) to the SwitchableRangedWeapon sheet,
# documenting the bot-only Semi/Burst fire-interval parameter, and tidies up
# a couple of border quirks that came along with the manual row insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Excel border-edge constants (xlEdgeLeft/Top/Bottom/Right) and line styles.
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlLineStyleNone = -4142
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. New row 13 content.
#    Values are written in this particular order so that the new shared
#    strings land at the same indices Excel produced (BotReload, then the
#    "人机 " category, then the description, then the comment).
# ---------------------------------------------------------------------
$ws.Range("B13").Value = "BotReload"
$ws.Range("A13").Value = "人机 "
$ws.Range("D13").Value = "Semi/Burst模式下机器人开火间隔"
$ws.Range("H13").Value = "仅对机器人手中的Semi/Burst模式的武器生效"
$ws.Range("C13").Value = "为一个浮点数"
$ws.Range("E13").Value = "/"
$ws.Range("F13").Value = "/"
$ws.Range("G13").Value = 0.1

# ---------------------------------------------------------------------
# 2. Formatting for the new row - centered like every other row, with a
#    thin box border around A13:H13 (same look as the row 12 box).
# ---------------------------------------------------------------------
$newRow = $ws.Range("A13:H13")
$newRow.HorizontalAlignment = $xlCenter
$newRow.VerticalAlignment = $xlCenter

$ws.Range("A13").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("A13").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("A13").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

$ws.Range("B13").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

$ws.Range("C13:G13").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("C13:G13").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

$ws.Range("H13").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("H13").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("H13").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

# ---------------------------------------------------------------------
# 3. Border clean-up left over from the manual row insertion:
#    - row 12's box now continues down into row 13, so its bottom edge
#      along column G opens up;
#    - C12 loses the stray left divider it inherited from column A;
#    - column B loses its divider against column A for the two
#      merged-category blocks (rows 2-5 and row 11).
# ---------------------------------------------------------------------
$ws.Range("G12").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone
$ws.Range("C12").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone

$ws.Range("B2").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone
$ws.Range("B3").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone
$ws.Range("B4").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone
$ws.Range("B5").Borders.Item($xlEdgeLeft).LineStyle = $xlLineStyleNone
$ws.Range("B11").Borders.Item($xlEdgeBottom).LineStyle = $xlLineStyleNone

# ---------------------------------------------------------------------
# 4. Restore the default view (scrolled back to column A, selection left
#    where the editing session ended up).
# ---------------------------------------------------------------------
$ws.Range("D30").Select()
